$wb = $excel.ActiveWorkbook

# --- 1. Update the Status text from "Ready for handoff" to "In Translation"
#     on every sheet where it appears (Overview!E2:F2, zh-cn!C2, de-de!C2). ---

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2. The status column(s) got narrower once the new text is shorter than
#     "Ready for handoff" (column autofit). Apply the resulting width. ---

$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
